$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at row 9, shifting existing rows 9-18 down to 10-19
$ws.Rows.Item(9).Insert()

# Copy the style (number format) of the date cell from the row below (now row 10) into the new row 9
$ws.Range("D10").Copy()
$ws.Range("D9").PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = $false

# Fill in the new record for row 9
$ws.Cells.Item(9, 1).Value = 7
$ws.Cells.Item(9, 2).Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Cells.Item(9, 3).Value = "Ñuble"
$ws.Cells.Item(9, 4).Value = 44771
$ws.Cells.Item(9, 5).Value = 16
$ws.Cells.Item(9, 6).Value = 100112037
$ws.Cells.Item(9, 7).Value = "Cebollín"
$ws.Cells.Item(9, 8).Value = "Sin especificar"
$ws.Cells.Item(9, 9).Value = "Primera"
$ws.Cells.Item(9, 10).Value = 150
$ws.Cells.Item(9, 11).Value = 8000
$ws.Cells.Item(9, 12).Value = 8000
$ws.Cells.Item(9, 13).Value = 8000
$ws.Cells.Item(9, 14).Value = "$/docena de atados"
$ws.Cells.Item(9, 15).Value = "Provincia de Diguillín"
$ws.Cells.Item(9, 16).Value = 2667
$ws.Cells.Item(9, 17).Value = 3
$ws.Cells.Item(9, 18).Value = "Hortaliza"
